$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the timestamp for the existing row (row 2) to the new run time
$ws.Range("A2").Value = "2025-09-14 12:40:05"

# Widen columns B and D slightly (23 -> 29, 28 -> 30 in saved width units).
# Excel's ColumnWidth property measured in characters is rounded/padded on
# save, so subtract 0.9 from the desired saved width to land on the exact
# target value once persisted.
$ws.Columns.Item(2).ColumnWidth = 28.1
$ws.Columns.Item(4).ColumnWidth = 29.1

# Append the new scraped listing as row 3
$ws.Range("A3").Value = "2025-09-14 12:40:05"
$ws.Range("B3").Value = "【医療関連】会員制サイト構築のパートナーを探しています"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5393406"
$ws.Range("G3").Value = 45
$ws.Range("H3").Value = "◇サイト"

# Add the hyperlink for the new URL cell, then restore the shared
# "Hyperlink" cell style so it matches the existing F2 cell (Add() alone
# creates a fresh style entry instead of reusing the existing one).
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5393406")
$ws.Range("F3").Style = "Hyperlink"
